# ICA02-Uke5-CircuitBreakers.xlsx — "Adding ica02 encoder / decoder + readme"
#
# The sheet gains a "Sannsynlighet" header above the existing "Per 100
# studenter" column (which itself gets relabeled "Per 100 studenter i %"),
# plus a new bottom "Gjennomsnitt" (average) summary row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row tweaks -----------------------------------------------
# New label above the existing "Per 100 studenter" header (row 1, col P).
$ws.Range("P1").Value = "Sannsynlighet"
# Existing header (row 2, col P) gets a more descriptive label.
$ws.Range("P2").Value = "Per 100 studenter i %"

# --- New "Gjennomsnitt" (average) summary row -------------------------
$ws.Range("A10").Value = "Gjennomsnitt"

$ws.Range("E10").Formula = "=AVERAGE(E3:E9)"
$ws.Range("G10").Formula = "=AVERAGE(G3:G9)"
$ws.Range("J10").Formula = "=AVERAGE(J3:J9)"
$ws.Range("L10").Formula = "=AVERAGE(L3:L9)"
$ws.Range("N10").Formula = "=AVERAGE(N3:N9)"
$ws.Range("P10").Formula = "=AVERAGE(P3:P9)"

# Bold the new average figures (not the "Gjennomsnitt" label itself).
$ws.Range("E10").Font.Bold = $true
$ws.Range("G10").Font.Bold = $true
$ws.Range("J10").Font.Bold = $true
$ws.Range("L10").Font.Bold = $true
$ws.Range("N10").Font.Bold = $true
$ws.Range("P10").Font.Bold = $true

# --- Page setup ---------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection ------------------------------------------------------
$ws.Range("G17").Select() | Out-Null
